$wb = $excel.ActiveWorkbook

# --- Customers sheet (sheet1) ---
$wsCustomers = $wb.Worksheets.Item("Customers")

# Row 2
$wsCustomers.Range("A2").Value = 1
$wsCustomers.Range("B2").Value = "A"
$wsCustomers.Range("C2").Value = "B"

# Row 3
$wsCustomers.Range("A3").Value = 2
$wsCustomers.Range("B3").Value = "C"
$wsCustomers.Range("C3").Value = "D"

# Email column (written after A-D columns so shared-string order matches)
$wsCustomers.Range("D2").Value = "ab@g.com"
$wsCustomers.Range("D3").Value = "cb@g.com"

# Hyperlinks on the email column
$wsCustomers.Hyperlinks.Add($wsCustomers.Range("D2"), "mailto:ab@g.com") | Out-Null
$wsCustomers.Hyperlinks.Add($wsCustomers.Range("D3"), "mailto:cb@g.com") | Out-Null

$wsCustomers.Range("E2").Select() | Out-Null

# --- Devices sheet (sheet2) ---
$wsDevices = $wb.Worksheets.Item("Devices")

# Update header row: insert device_name before device_measurements
$wsDevices.Range("C1").Value = "device_name"
$wsDevices.Range("D1").Value = "device_measurements"

# Device rows
$wsDevices.Range("A2").Value = 1
$wsDevices.Range("B2").Value = 1
$wsDevices.Range("C2").Value = "ESP1"

$wsDevices.Range("A3").Value = 2
$wsDevices.Range("B3").Value = 1
$wsDevices.Range("C3").Value = "ESP3"

$wsDevices.Range("A4").Value = 3
$wsDevices.Range("B4").Value = 1
$wsDevices.Range("C4").Value = "Ardunio"

$wsDevices.Range("A5").Value = 4
$wsDevices.Range("B5").Value = 2
$wsDevices.Range("C5").Value = "Telus"

$wsDevices.Range("A6").Value = 5
$wsDevices.Range("B6").Value = 2
$wsDevices.Range("C6").Value = "Rogers"

$wsDevices.Range("A7").Value = 6
$wsDevices.Range("B7").Value = 2
$wsDevices.Range("C7").Value = "Alpha"

$wsDevices.Range("A8").Value = 7
$wsDevices.Range("B8").Value = 2
$wsDevices.Range("C8").Value = "Beta"

$wsDevices.Range("A9").Value = 8
$wsDevices.Range("B9").Value = 2
$wsDevices.Range("C9").Value = "Gamma"

$wsDevices.Columns.Item(2).AutoFit() | Out-Null
$wsDevices.Columns.Item(3).AutoFit() | Out-Null

$wsDevices.Range("C10").Select() | Out-Null
